# GDE-9501: Add TL Base Rate for PT Health
#
# Inserts a new "BaseRate_Fields" test-case row (rowid 7) that points at the
# 27-Dec-2019 base-rate template/GS-file path combo but references the new
# GS input file FINASTRA_CCB_BASERATE_SY_GROUP1_20191125.csv. The two rows
# that used to be rowid 7/8 shift down to rows 9/10 and are renumbered to
# rowid 8/9 respectively.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BaseRate_Fields")
$ws2 = $wb.Worksheets.Item("FXRates_Fields")

# Push the old rows 8 & 9 down to make room for the new row 8 (this also
# carries the existing number formats/styles and shifts the existing
# hyperlink ranges down by one row, same as Excel's native "Insert Rows").
$ws1.Rows.Item(8).Insert()

# Populate the new row with the same TemplateFilePath pair used by the
# 27-Dec-2019 base rate rows (rows 2-7).
$ws1.Range("A8").Value = "7"
$ws1.Range("B8").Value = "01_TL_Base_Rates_27122019"

# TemplateFilePath/InputFilePath (columns C & D) carry a hyperlink to the
# shared GS-file network location, same as every other row in the sheet.
# Add the hyperlinks first (this seeds the cell text too), then overwrite
# the cell text with the real path - this matches how the rest of the
# sheet's hyperlinks ended up with a display name independent of the
# cell's actual text.
$ws1.Hyperlinks.Add($ws1.Range("D8"), "file:///\\DataSet\TL_DataSet\BaseRates_GSFile\", "", "", "\\DataSet\\TL_DataSet\\BaseRates_GSFile\\") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("C8"), "file:///\\DataSet\TL_DataSet\BaseRates_GSFile\", "", "", "\\DataSet\\TL_DataSet\\BaseRates_GSFile\\") | Out-Null

$ws1.Range("C8").Value = "\DataSet\NewUATDeals_DataSet\Transformation_Layer\TL_Base_Rate\"
$ws1.Range("D8").Value = "\DataSet\NewUATDeals_DataSet\Transformation_Layer\TL_Base_Rate\BaseRates_Files_27DEC2019\"

# New GS input file for the newly inserted test case.
$ws1.Range("E8").Value = "FINASTRA_CCB_BASERATE_SY_GROUP1_20191125.csv"

# Renumber the two rows that were pushed down.
$ws1.Range("A9").Value = "8"
$ws1.Range("A10").Value = "9"

# Writing .Value resets a cell's pre-existing style (quote-prefixed text for
# column A, plain text for C/D once the hyperlink auto-formatting kicks in)
# so re-apply the correct format from an untouched neighbour cell.
$ws1.Range("A7").Copy()
$ws1.Range("A8").PasteSpecial(-4122)
$ws1.Range("A9").PasteSpecial(-4122)
$ws1.Range("A10").PasteSpecial(-4122)

$ws1.Range("C7").Copy()
$ws1.Range("C8").PasteSpecial(-4122)
$ws1.Range("D7").Copy()
$ws1.Range("D8").PasteSpecial(-4122)

# Leave the selection where the edit ended up: E9 (old row 8's GS-file cell,
# now shifted to row 9) on the base-rate sheet, and restore the FX sheet's
# own last-used selection without leaving it as the active tab.
$ws2.Range("B34").Select()
$ws1.Range("E9").Select()
